$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.226.59'
$ws.Range("E2").Value = '  +2.95%  '
$ws.Range("D3").Value = '2.072.50'
$ws.Range("E3").Value = '  +2.72%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.43'
$ws.Range("E5").Value = '  +2.11%  '
$ws.Range("E6").Value = '  +1.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.16'
$ws.Range("E7").Value = '  +5.86%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +2.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0809'
$ws.Range("E10").Value = '  +2.66%  '
$ws.Range("D12").Value = '2.378.62'
$ws.Range("E12").Value = '  +2.59%  '
$ws.Range("E13").Value = '  +2.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.83'
$ws.Range("E14").Value = '  +1.86%  '
$ws.Range("E15").Value = '  +1.43%  '
$ws.Range("E16").Value = '  +2.48%  '
$ws.Range("D17").Value = '2.065.45'
$ws.Range("E17").Value = '  +3.48%  '
$ws.Range("D18").Value = '38.112.36'
$ws.Range("E18").Value = '  +2.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.24'
$ws.Range("E19").Value = '  +1.25%  '
$ws.Range("E20").Value = '  +1.69%  '
$ws.Range("D21").Value = '0.0₃0834'
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.38'
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  +2.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.34'
$ws.Range("E26").Value = '  +1.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.11'
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("E28").Value = '  +7.07%  '
$ws.Range("E29").Value = '  +2.19%  '
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("E31").Value = '  +1.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.58'
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.61'
$ws.Range("E33").Value = '  +4.19%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0617'
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("E35").Value = '  +7.93%  '
$ws.Range("E36").Value = '  +1.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.05'
$ws.Range("E37").Value = '  +11.91%  '
$ws.Range("E38").Value = '  +4.83%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '98.63'
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0220'
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("D42").Value = '1.484.46'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0952'
$ws.Range("E43").Value = '  +2.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.85'
$ws.Range("E44").Value = '  +1.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.87'
$ws.Range("E45").Value = '  +3.70%  '
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.05'
$ws.Range("E47").Value = '  +15.28%  '
$ws.Range("E48").Value = '  +2.08%  '
$ws.Range("E49").Value = '  +1.96%  '
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("D51").Value = '2.262.86'
$ws.Range("E51").Value = '  +2.43%  '
